# Automatische test-sync: 2025-08-05 18:24:50
$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: add new row 30 (Testmail #9) ---
$logs.Range("A30").Value = "Hoi, hebben jullie al iets gehoord?"
$logs.Range("B30").Value = "mailmind.test@zohomail.eu"
$logs.Range("C30").Value = "Testmail #9: Hoi, hebben jullie al iets gehoord?"
$logs.Range("D30").Value = "Opvolging / Status"
$logs.Range("E30").Value = "Dank voor je bericht. We hebben je eerdere e-mail ontvangen en doorgestuurd naar klantenservice@bedrijf.nl."
$logs.Range("F30").Value = "2025-08-05 18:23:51"
$logs.Range("G30").Value = "Ja"
$logs.Range("H30").Value = "Ja"
$logs.Range("I30").Value = "Nee"
$logs.Range("J30").Value = "Nee"

# --- Logs sheet: extend conditional formatting ranges to include row 30 ---
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "29")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "30")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard sheet: swap row 3 / row 4 category labels ---
$dash.Range("A3").Value = "Klantenservice / Contact"
$dash.Range("A4").Value = "Inkoop / Bestellingen"

# --- Dashboard sheet: add new row 7 (Opvolging / Status) ---
$dash.Range("A7").Value = "Opvolging / Status"
$dash.Range("B7").Value = 1

# --- Chart1: extend category/value series ranges from row 6 to row 7 ---
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.Values = "='Dashboard'!`$B`$2:`$B`$7"
$series.XValues = "='Dashboard'!`$A`$2:`$A`$7"
